$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.55%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.02%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.700"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'7.81%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08005"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.25%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.12%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.492"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.08%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.973"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9255"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.10%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1240"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-7.46%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-0.07%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'8.720"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'21.72%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09166"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.60%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.03640"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.97%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.1050"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'9.62%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-2.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006136"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.79%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.351"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.71%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3474"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.36%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'3.82%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004648"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'8.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001130"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-5.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02492"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.43%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05330"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.44%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007465"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.28%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009642"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.04%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.51%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-2.50%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01034"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.30%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006714"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.80%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002971"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-11.18%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-4.59%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
